$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row for ID "QSY_B_072" (row 14). This shifts the
# rows below it up by one: QSY_B_189 becomes row 14 and the trailing
# blank row becomes row 15.
$ws.Rows(14).Delete()

# Match the author's final cell selection recorded in the saved file.
$ws.Range("M14").Select()
